$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row just above the "unit_costs.tsv" block (old row 16),
# pushing it (and everything after) down by one. Rows 13-15 ("cia_factbook_note",
# "intermediate_area_name", "intermediate_area_count") are left untouched by
# this step, and the freshly-inserted row 16 inherits formatting (style indices)
# from the row above.
$ws.Rows("16:16").Insert()
$ws.Rows("16:16").RowHeight = 90
# The insert synthesizes an (empty) D16 cell that the target layout doesn't
# have -- drop it so the row only carries the A/B/C cells it needs.
$ws.Range("D16").Clear()

# Now re-populate rows 13-16 with the correct (reordered + new) content.
# Row 13: intermediate_area_name (used to be row 14)
$ws.Range("A13").Value = "country.tsv"
$ws.Range("B13").Value = "intermediate_area_name"
$ws.Range("C13").Value = "The name of intermediate areas identified by the CIA World Factbook data. Note that while these data can be used to roughly estimate the number of administrative regions per country; this information is best supplemented by local expertise and information on the administrative organization of IHR-related activities. "
$ws.Range("D13").Value = "From the CIA World Factbook: based on ""designatory terms, and first-order administrative divisions as approved by the US Board on Geographic Names (BGN)"""

# Row 14: intermediate_area_count (used to be row 15)
$ws.Range("A14").Value = "country.tsv"
$ws.Range("B14").Value = "intermediate_area_count"
$ws.Range("C14").Value = "The number of intermediate areas identified by the CIA World Factbook data. Note that while these data can be used to roughly estimate the number of administrative regions per country; this information is best supplemented by local expertise and information on the administrative organization of IHR-related activities. "
$ws.Range("D14").Value = "From the CIA World Factbook: based on ""designatory terms, and first-order administrative divisions as approved by the US Board on Geographic Names (BGN)"""

# Row 15: brand-new field documenting the source of the intermediate area data
$ws.Range("A15").Value = "country.tsv"
$ws.Range("B15").Value = "intermediate_area_reference"
$ws.Range("C15").Value = "Metadata documenting the source (and date) based on which the intermediate area data were determined"
$ws.Range("D15").Value = "Currently all intermediate data come from the CIA World Factbook (2022)"

# Row 16: cia_factbook_note (used to be row 13), no notes column for this field
$ws.Range("A16").Value = "country.tsv"
$ws.Range("B16").Value = "cia_factbook_note"
$ws.Range("C16").Value = "Any notes included in the data export from the CIA World Factbook and/or any notes made by the research team during manual extraction of administrative area data."

# Match the saved view state: scrolled so row 13 is at the top, with C15 selected.
$ws.Range("C15").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
